$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the listed rows to reflect repulled data / mean calculation
$ws.Range("F5").Value = -3
$ws.Range("F11").Value = 0
$ws.Range("F20").Value = -4
$ws.Range("F24").Value = 1
$ws.Range("F30").Value = 2
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 3
$ws.Range("F41").Value = -2
$ws.Range("F43").Value = 2
$ws.Range("F45").Value = -1
$ws.Range("F46").Value = 0
$ws.Range("F50").Value = -2
$ws.Range("F62").Value = -1
$ws.Range("F63").Value = 2
$ws.Range("F66").Value = 0
$ws.Range("F68").Value = 1
